# Generate Report for Handback
# Fill in the "Latest Target File" / "Latest Handback File" / "Latest
# Handback DateTime" / "Error Detail" columns for the 6b2059c8... row on
# both the zh-cn and de-de sheets, now that a (stale) handback arrived.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # RGB(0x64,0x95,0xED) == FF6495ED, matches the
                             # workbook's existing "HyperLink" cell style

$mismatchMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54dbe5152e539a9fdf9034cad85e4d9a6ae97b92/e2e/6b2059c8-915f-4fd0-b7a9-2ecd629fee9d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a73fd263688ad87470f3af84e4eca9960ffd0da5/e2e/6b2059c8-915f-4fd0-b7a9-2ecd629fee9d.md."

# --- zh-cn sheet, row 7 (6b2059c8-915f-4fd0-b7a9-2ecd629fee9d) ---
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Add($ws.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a73fd263688ad87470f3af84e4eca9960ffd0da5/e2e/6b2059c8-915f-4fd0-b7a9-2ecd629fee9d.md", "", "", "6b2059c8-915f-4fd0-b7a9-2ecd629fee9d.md")
$ws.Range("I7").Font.Underline = $true
$ws.Range("I7").Font.Color = $hyperlinkColor

$ws.Range("J7").Value = "6b2059c8-915f-4fd0-b7a9-2ecd629fee9d.3f2eff76d3bb12b08c5c0795bc0ace214ddeff7d.zh-cn.xlf"
$ws.Range("K7").Value = "2016-09-01 23:00:30"
$ws.Range("P7").Value = $mismatchMessage

# --- de-de sheet, row 7 (6b2059c8-915f-4fd0-b7a9-2ecd629fee9d) ---
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Hyperlinks.Add($ws2.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a73fd263688ad87470f3af84e4eca9960ffd0da5/e2e/6b2059c8-915f-4fd0-b7a9-2ecd629fee9d.md", "", "", "6b2059c8-915f-4fd0-b7a9-2ecd629fee9d.md")
$ws2.Range("I7").Font.Underline = $true
$ws2.Range("I7").Font.Color = $hyperlinkColor

$ws2.Range("J7").Value = "6b2059c8-915f-4fd0-b7a9-2ecd629fee9d.3f2eff76d3bb12b08c5c0795bc0ace214ddeff7d.de-de.xlf"
$ws2.Range("K7").Value = "2016-09-01 23:00:38"
$ws2.Range("P7").Value = $mismatchMessage

Write-Host "Report generated for handback row."
